{"js": "const replacements = [\n  [\"845\u00d74=3380\", \"606\u00d79=5454\"],\n  [\"296\u00d79=2664\", \"625\u00d79=5625\"],\n  [\"435\u00d74=1740\", \"789\u00d75=3945\"],\n  [\"143\u00d75=715\", \"424\u00d73=1272\"],\n  [\"918\u00d72=1836\", \"915\u00d72=1830\"],\n  [\"447\u00d72=894\", \"311\u00d72=622\"],\n  [\"684\u00d74=2736\", \"544\u00d76=3264\"],\n  [\"284\u00d73=852\", \"313\u00d75=1565\"],\n  [\"614\u00d78=4912\", \"989\u00d78=7912\"],\n  [\"247\u00d72=494\", \"866\u00d77=6062\"],\n  [\"524\u00d77=3668\", \"665\u00d72=1330\"],\n  [\"789\u00d73=2367\", \"193\u00d78=1544\"],\n  [\"413\u00d75=2065\", \"695\u00d79=6255\"],\n  [\"467\u00d76=2802\", \"278\u00d77=1946\"],\n  [\"548\u00d75=2740\", \"895\u00d78=7160\"],\n  [\"486\u00d76=2916\", \"506\u00d77=3542\"],\n  [\"310\u00d76=1860\", \"764\u00d73=2292\"],\n  [\"961\u00d73=2883\", \"477\u00d72=954\"],\n  [\"270\u00d76=1620\", \"605\u00d74=2420\"],\n  [\"298\u00d72=596\", \"766\u00d77=5362\"],\n  [\"904\u00d72=1808\", \"949\u00d73=2847\"],\n  [\"159\u00d77=1113\", \"399\u00d79=3591\"],\n  [\"543\u00d79=4887\", \"693\u00d74=2772\"],\n  [\"436\u00d78=3488\", \"462\u00d78=3696\"],\n  [\"741\u00d79=6669\", \"816\u00d79=7344\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"845\u00d74=3380\", \"606\u00d79=5454\"),\n    @(\"296\u00d79=2664\", \"625\u00d79=5625\"),\n    @(\"435\u00d74=1740\", \"789\u00d75=3945\"),\n    @(\"143\u00d75=715\", \"424\u00d73=1272\"),\n    @(\"918\u00d72=1836\", \"915\u00d72=1830\"),\n    @(\"447\u00d72=894\", \"311\u00d72=622\"),\n    @(\"684\u00d74=2736\", \"544\u00d76=3264\"),\n    @(\"284\u00d73=852\", \"313\u00d75=1565\"),\n    @(\"614\u00d78=4912\", \"989\u00d78=7912\"),\n    @(\"247\u00d72=494\", \"866\u00d77=6062\"),\n    @(\"524\u00d77=3668\", \"665\u00d72=1330\"),\n    @(\"789\u00d73=2367\", \"193\u00d78=1544\"),\n    @(\"413\u00d75=2065\", \"695\u00d79=6255\"),\n    @(\"467\u00d76=2802\", \"278\u00d77=1946\"),\n    @(\"548\u00d75=2740\", \"895\u00d78=7160\"),\n    @(\"486\u00d76=2916\", \"506\u00d77=3542\"),\n    @(\"310\u00d76=1860\", \"764\u00d73=2292\"),\n    @(\"961\u00d73=2883\", \"477\u00d72=954\"),\n    @(\"270\u00d76=1620\", \"605\u00d74=2420\"),\n    @(\"298\u00d72=596\", \"766\u00d77=5362\"),\n    @(\"904\u00d72=1808\", \"949\u00d73=2847\"),\n    @(\"159\u00d77=1113\", \"399\u00d79=3591\"),\n    @(\"543\u00d79=4887\", \"693\u00d74=2772\"),\n    @(\"436\u00d78=3488\", \"462\u00d78=3696\"),\n    @(\"741\u00d79=6669\", \"816\u00d79=7344\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $found = $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Replacement text not found: $oldText\"\n    }\n}\n"}
